# Upstream feed update (bilibili 会员购 "苏州" listings): the 2024-04-13
# entry "苏州·绘时国乙1.0-秩序之外" was removed from the source data. That
# removes one row from every per-event sheet that listed it ("展览" and
# "全部类型" - the "演出"/"本地生活" sheets never contained this event).
#
# Column A is a plain sequential row index (0, 1, 2, ...) that is
# regenerated fresh for each row position and is NOT part of the shifted
# event data, so after removing the row we rewrite it back to a simple
# 0-based sequence for every remaining data row.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the row for "苏州·绘时国乙1.0-秩序之外" (2024-04-13); everything
    # below shifts up by one row, and the used range shrinks by one row.
    $ws.Rows.Item(3).Delete()

    # Restore column A as a plain 0-based row index (row 1 is the header
    # row, value 0; row 2 is the first data row, value 1; etc).
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}
